# Add ability to support multiple entries per form.
# The sheet originally stores a single form "entry" per response row in
# columns J:S (Entry Date + 9 question-option columns). This change adds
# two more copies of that same 10-column block (T:AC and AD:AM) so each
# row can hold up to three form entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 headers for the two new blocks (9 option columns each; the
#      "Entry Date" header itself is not repeated for the new blocks) ----
$optionHeaders = @("q1_opt1","q1_opt2","q1_opt3","q2_opt1","q2_opt2","q2_opt3","q3_opt1","q3_opt2","q3_opt3")

# Block 1 headers -> columns T..AB (20..28)
for ($i = 0; $i -lt $optionHeaders.Length; $i++) {
    $ws.Cells.Item(1, 20 + $i).Value = $optionHeaders[$i]
}
# Block 2 headers -> columns AC..AK (29..37)
for ($i = 0; $i -lt $optionHeaders.Length; $i++) {
    $ws.Cells.Item(1, 29 + $i).Value = $optionHeaders[$i]
}

# ---- Row 3 data: second entry (T:AC) + third entry (AD:AM) ----
$row3Block1 = @(1, 0, 0, 0, 1, 0, 2, 3, 5)
$row3Block2 = @(1, 0, 0, 1, 0, 0, 1, 1, 1)

for ($i = 0; $i -lt $row3Block1.Length; $i++) {
    $ws.Cells.Item(3, 21 + $i).Value = $row3Block1[$i]
}
for ($i = 0; $i -lt $row3Block2.Length; $i++) {
    $ws.Cells.Item(3, 31 + $i).Value = $row3Block2[$i]
}

# ---- Row 4 data: second entry (T:AC) + third entry (AD:AM) ----
$row4Block1 = @(1, 0, 0, 0, 1, 0, 2, 3, 4)
$row4Block2 = @(0, 0, 1, 0, 0, 1, 3, 3, 3)

for ($i = 0; $i -lt $row4Block1.Length; $i++) {
    $ws.Cells.Item(4, 21 + $i).Value = $row4Block1[$i]
}
for ($i = 0; $i -lt $row4Block2.Length; $i++) {
    $ws.Cells.Item(4, 31 + $i).Value = $row4Block2[$i]
}

# ---- The new "Entry Date" cells (T3, AD3, T4, AD4) use the same date
#      number format as the existing Entry Date column (J). Copy the
#      format from J3 (which already carries the date style) onto each new
#      date cell - this reuses the workbook's existing date style instead
#      of creating a new one - then set its value. ----
$dateCells = @("T3", "AD3", "T4", "AD4")
$dateValues = @(42707.97142585126, 42708.016006341655, 42707.97153375884, 42708.01611987991)
for ($i = 0; $i -lt $dateCells.Length; $i++) {
    $ws.Range("J3").Copy()
    $ws.Range($dateCells[$i]).PasteSpecial(-4122)
    $ws.Range($dateCells[$i]).Value = $dateValues[$i]
}

# ---- Column widths for the newly introduced columns (T..AM) ----
$ws.Columns("T:AM").AutoFit()

$excel.CutCopyMode = $false
